# Update the August (month 8) daily-sales figures for days 4, 5 and 6,
# which now reflect revised totals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 30338.81
$ws.Range("B4").Value = 16137.39
$ws.Range("B5").Value = 19849.77

# Two new daily records (Aug 7 and Aug 8) were added to the source data.
# They belong right after the existing August rows (row 5) and before the
# July rows (old row 6), so insert two blank rows at row 6 and fill them in.
$ws.Range("A6:A7").EntireRow.Insert()

$ws.Range("A6").Value = 7
$ws.Range("B6").Value = 22763.35
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = 2025
$ws.Range("E6").Value = "08/2025"

$ws.Range("A7").Value = 8
$ws.Range("B7").Value = 13507.96
$ws.Range("C7").Value = 8
$ws.Range("D7").Value = 2025
$ws.Range("E7").Value = "08/2025"
